# Update cryptocurrency price/volume symbol data on Sheet1 to reflect the
# latest scrape performed by the GitHub Actions workflow.
#
# The "Price" column (D) and some "Volume(1h)" column (E) cells are stored
# as text (not numbers), so we force a text number format before assigning
# the new values to make sure Excel keeps them as strings instead of
# re-interpreting them as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2"  "274.87"
Set-TextValue "D3"  "23.08"
Set-TextValue "D5"  "0.06276"
Set-TextValue "D6"  "3.656"
Set-TextValue "D7"  "6.678"
Set-TextValue "D8"  "1.392"
Set-TextValue "D9"  "0.8345"
Set-TextValue "D10" "0.01388"
Set-TextValue "D11" "0.1623"
Set-TextValue "D12" "0.08269"
Set-TextValue "D13" "0.03427"
Set-TextValue "D14" "0.03115"
Set-TextValue "D15" "0.09296"
Set-TextValue "D16" "3.855"
Set-TextValue "D17" "0.001652"
Set-TextValue "D18" "0.04774"
Set-TextValue "D19" "0.006393"

Set-TextValue "E20" "19HotbitTokenHTBWorstin24h"

Set-TextValue "D21" "0.001087"
Set-TextValue "D23" "3.714"
Set-TextValue "D25" "0.3345"
Set-TextValue "D26" "0.1256"
Set-TextValue "D27" "0.0002681"
Set-TextValue "D40" "0.04712"
Set-TextValue "D41" "0.007042"
Set-TextValue "D42" "0.1163"
Set-TextValue "D43" "0.003351"
Set-TextValue "D44" "0.01211"
Set-TextValue "D45" "0.00006269"
Set-TextValue "D48" "0.7968"

Set-TextValue "E48" "47CoinbaseStockTokenCOIN"

Set-TextValue "D49" "0.006605"
Set-TextValue "D50" "0.00002301"
Set-TextValue "D51" "0.01241"
